# Applies the coin price/volume refresh captured in the commit message:
# "Updated symbol list on Wed Feb  1 14:56:16 UTC 2023 with GitHub Actions".
# Only the Price (column D) and Volume(1h) (column E) cells change, for a
# subset of the listed coins (rows 2-27 and 39-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price / Volume(1h) text for every cell that actually changed.
# (Both columns hold plain text such as "309.92" / "-0.60%", so values
# are kept as Text below instead of being auto-converted by Excel into
# Number / Percentage values.)
$newValues = [ordered]@{
    "D2" = "309.83"
    "E2" = "-0.68%"
    "D3" = "37.23"
    "E3" = "-2.09%"
    "D4" = "5.131"
    "E4" = "0.06%"
    "D5" = "0.07846"
    "E5" = "-1.10%"
    "D6" = "8.277"
    "E6" = "0.43%"
    "D7" = "1.880"
    "E7" = "-1.73%"
    "E8" = "4.32%"
    "D9" = "0.9249"
    "E9" = "-0.17%"
    "D10" = "0.1178"
    "E10" = "-2.15%"
    "D11" = "0.1896"
    "E11" = "-0.89%"
    "D12" = "0.08879"
    "E12" = "-3.66%"
    "D13" = "0.03312"
    "E13" = "-2.46%"
    "D14" = "0.09612"
    "E14" = "-0.14%"
    "E15" = "0.31%"
    "D16" = "0.006196"
    "E16" = "5.91%"
    "D17" = "3.393"
    "E17" = "-3.97%"
    "D18" = "4.400"
    "E18" = "-0.02%"
    "E19" = "0.41%"
    "D20" = "6.387"
    "E20" = "21.43%"
    "D21" = "0.1292"
    "E21" = "0.70%"
    "D22" = "0.2407"
    "E22" = "-6.96%"
    "D23" = "0.04342"
    "E23" = "-0.56%"
    "E24" = "-3.87%"
    "D25" = "0.004285"
    "E25" = "0.12%"
    "D26" = "0.0001400"
    "D27" = "0.0002901"
    "D39" = "0.02160"
    "E39" = "2.56%"
    "D40" = "0.05008"
    "E40" = "-1.52%"
    "D41" = "0.007585"
    "E41" = "-0.66%"
    "D42" = "0.1357"
    "E42" = "0.15%"
    "D43" = "0.008484"
    "E43" = "-7.02%"
    "D44" = "0.002011"
    "E44" = "-1.09%"
    "D45" = "0.007982"
    "E45" = "-7.67%"
    "D46" = "0.00006573"
    "E46" = "-1.54%"
    "D47" = "0.00000000750"
    "E47" = "0.12%"
    "D48" = "0.003294"
    "E48" = "13.81%"
    "D49" = "0.001443"
    "E49" = "20.39%"
    "D50" = "0.00002101"
    "E50" = "0.12%"
    "D51" = "0.0002001"
    "E51" = "0.12%"
}

foreach ($ref in $newValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$ref]
}
